$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "February "
$ws.Range("D2").Value = 83

$ws.Range("A3").Value = "Revolution "
$ws.Range("B3").Value = 152
$ws.Range("C3").Value = 669
$ws.Range("D3").Value = 98

$ws.Range("A4").Value = ". "
$ws.Range("B4").Value = 250
$ws.Range("C4").Value = 669
$ws.Range("D4").Value = 10

$ws.Range("A5").Value = "The "
$ws.Range("B5").Value = 260
$ws.Range("C5").Value = 669
$ws.Range("D5").Value = 37

$ws.Range("A6").Value = "second "
$ws.Range("B6").Value = 297
$ws.Range("C6").Value = 669
$ws.Range("D6").Value = 64
$ws.Range("E6").Value = 23

$ws.Range("A7").Value = "was "
$ws.Range("B7").Value = 361
$ws.Range("C7").Value = 669
$ws.Range("D7").Value = 37
$ws.Range("E7").Value = 23

$ws.Range("A8").Value = "the "
$ws.Range("B8").Value = 398
$ws.Range("C8").Value = 669
$ws.Range("D8").Value = 32
$ws.Range("E8").Value = 23

$ws.Range("A9").Value = "October "
$ws.Range("B9").Value = 430
$ws.Range("C9").Value = 669
$ws.Range("D9").Value = 73
$ws.Range("E9").Value = 23

$ws.Range("A10").Value = "Revolution."
$ws.Range("B10").Value = 503
$ws.Range("C10").Value = 669
$ws.Range("D10").Value = 99
$ws.Range("E10").Value = 23

$ws.Range("A11").Value = "czarist "
$ws.Range("B11").Value = 497
$ws.Range("C11").Value = 710.4
$ws.Range("D11").Value = 60
$ws.Range("E11").Value = 23

$ws.Range("A12").Value = "government."
$ws.Range("B12").Value = 557
$ws.Range("C12").Value = 710.4
$ws.Range("D12").Value = 108
$ws.Range("E12").Value = 23

$ws.Range("A13").Value = "Russian "
$ws.Range("B13").Value = 816
$ws.Range("C13").Value = 1324.4
$ws.Range("D13").Value = 73
$ws.Range("E13").Value = 23

$ws.Range("A14").Value = "Civil "
$ws.Range("B14").Value = 889
$ws.Range("C14").Value = 1324.4
$ws.Range("D14").Value = 44
$ws.Range("E14").Value = 23

$ws.Range("A15").Value = "War. "
$ws.Range("B15").Value = 152
$ws.Range("C15").Value = 1349.4
$ws.Range("D15").Value = 47
$ws.Range("E15").Value = 23

$ws.Range("A16").Value = "Soviet "
$ws.Range("B16").Value = 493
$ws.Range("C16").Value = 1349.4
$ws.Range("D16").Value = 57
$ws.Range("E16").Value = 23

$ws.Range("A17").Value = "Union "
$ws.Range("B17").Value = 550
$ws.Range("C17").Value = 1349.4
$ws.Range("D17").Value = 58
$ws.Range("E17").Value = 23
